$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header labels (plain text, default style)
$ws.Range("A1").Value = "dist_hybrid_rew"
$ws.Range("B1").Value = "dist_hybrid_rew"
$ws.Range("C1").Value = "dist_hybrid_rew"

# Row 2
$ws.Range("A2").Value = "'0.9444444444444444"
$ws.Range("B2").Value = "'0.9444444444444444"
$ws.Range("C2").Value = "'0.3333333333333333"

# Row 3
$ws.Range("A3").Value = "'0.6304347826086957"
$ws.Range("B3").Value = "'0.8152173913043478"
$ws.Range("C3").Value = 0.5

# Row 4
$ws.Range("A4").Value = "'0.7982456140350878"
$ws.Range("B4").Value = "'0.7982456140350878"
$ws.Range("C4").Value = "'0.16666666666666666"

# Row 5
$ws.Range("A5").Value = "'0.7530864197530864"
$ws.Range("B5").Value = "'0.7530864197530864"
$ws.Range("C5").Value = "'0.16666666666666666"

# Row 6
$ws.Range("A6").Value = "'0.7521367521367521"
$ws.Range("B6").Value = "'0.7521367521367521"
$ws.Range("C6").Value = "'0.16666666666666666"

# Row 7
$ws.Range("A7").Value = "'0.5053763440860215"
$ws.Range("B7").Value = "'0.5483870967741935"
$ws.Range("C7").Value = "'0.6666666666666666"

# Row 8
$ws.Range("A8").Value = "'0.5087719298245614"
$ws.Range("B8").Value = "'0.8157894736842105"
$ws.Range("C8").Value = 0.5

# Row 9
$ws.Range("A9").Value = "'0.6555555555555556"
$ws.Range("B9").Value = "'0.6555555555555556"
$ws.Range("C9").Value = "'0.16666666666666666"

# Row 10
$ws.Range("A10").Value = "'0.8365384615384616"
$ws.Range("B10").Value = "'0.8365384615384616"
$ws.Range("C10").Value = 0.0

# Row 11
$ws.Range("A11").Value = "'0.8387096774193549"
$ws.Range("B11").Value = "'0.8924731182795699"
$ws.Range("C11").Value = 0.5

# Row 12
$ws.Range("A12").Value = "'0.5701754385964912"
$ws.Range("B12").Value = "'0.6140350877192983"
$ws.Range("C12").Value = "'0.3333333333333333"

# Row 13
$ws.Range("A13").Value = "'0.7549019607843137"
$ws.Range("B13").Value = "'0.7941176470588235"
$ws.Range("C13").Value = "'0.3333333333333333"

# Row 14
$ws.Range("A14").Value = "'0.5824175824175825"
$ws.Range("B14").Value = "'0.5824175824175825"
$ws.Range("C14").Value = 0.5

# Row 15
$ws.Range("A15").Value = "'0.7476635514018691"
$ws.Range("B15").Value = "'0.7476635514018691"
$ws.Range("C15").Value = "'0.3333333333333333"

# Match the final selected cell recorded in the workbook (I5)
[void]$ws.Range("I5").Select()
